$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.572.32'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '2.666.57'
$ws.Range('E3').Value = '  -1.98%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '600.59'
$ws.Range('E5').Value = '  -1.89%  '
$ws.Range('D6').Value = '156.95'
$ws.Range('E6').Value = '  -1.38%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '0.621'
$ws.Range('E8').Value = '  +5.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.130'
$ws.Range('E9').Value = '  +2.84%  '
$ws.Range('D10').Value = '0.403'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('D11').Value = '5.84'
$ws.Range('E11').Value = '  -4.11%  '
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('D13').Value = '29.26'
$ws.Range('E13').Value = '  -3.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000198'
$ws.Range('E14').Value = '  -5.05%  '
$ws.Range('D15').Value = '3.146.42'
$ws.Range('E15').Value = '  -1.89%  '
$ws.Range('D16').Value = '65.450.78'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').Value = '2.674.18'
$ws.Range('E17').Value = '  -2.05%  '
$ws.Range('D18').Value = '12.77'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').Value = '4.79'
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').Value = '7.59'
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('D21').Value = '349.58'
$ws.Range('E21').Value = '  -3.68%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').Value = '69.72'
$ws.Range('E23').Value = '  -0.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000110'
$ws.Range('E24').Value = '  +2.48%  '
$ws.Range('D25').Value = '9.65'
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('E26').Value = '  -1.91%  '
$ws.Range('E27').Value = '  -3.22%  '
$ws.Range('E28').Value = '  -8.16%  '
$ws.Range('D29').Value = '8.08'
$ws.Range('E29').Value = '  -3.31%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  -3.16%  '
$ws.Range('D32').Value = '530.42'
$ws.Range('E32').Value = '  -1.94%  '
$ws.Range('D33').Value = '1.76'
$ws.Range('E33').Value = '  -3.47%  '
$ws.Range('D34').Value = '5.49'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('D35').Value = '6.44'
$ws.Range('E35').Value = '  -3.51%  '
$ws.Range('D36').Value = '0.422'
$ws.Range('E36').Value = '  -3.47%  '
$ws.Range('D37').Value = '20.44'
$ws.Range('E37').Value = '  -2.36%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '158.40'
$ws.Range('E39').Value = '  -3.12%  '
$ws.Range('D40').Value = '1.93'
$ws.Range('E40').Value = '  -4.05%  '
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('D42').Value = '42.83'
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').Value = '164.25'
$ws.Range('E43').Value = '  -4.29%  '
$ws.Range('D44').Value = '4.14'
$ws.Range('E44').Value = '  -1.41%  '
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('D46').Value = '0.0607'
$ws.Range('E46').Value = '  -2.04%  '
$ws.Range('D47').Value = '22.86'
$ws.Range('E47').Value = '  -3.72%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.640'
$ws.Range('E48').Value = '  -3.71%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '0.0258'
$ws.Range('E49').Value = '  -3.46%  '
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.10'
$ws.Range('E51').Value = '  -2.46%  '
